$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(10, 11, 428.3635523572899, 0.08611118000000001, 0.08611118000000001, 0.17222236, 642.5453285359349, 0.4, 0.45, 0.12, 0.4),
    @(11, 12, 498.7069374946091, 0.10000008, 0.10000008, 0.20000016, 748.0604062419137, 0.34, 0.42, 0.14, 0.3),
    @(12, 13, 508.3495138842192, 0.11944454, 0.11944454, 0.23888908, 762.5242708263288, 0.35, 0.48, 0.12, 0.36),
    @(13, 14, 425.5577207029365, 0.10277786, 0.10277786, 0.20555572, 638.3365810544047, 0.35, 0.45, 0.14, 0.38),
    @(14, 15, 334.7633794933903, 0.0416667, 0.0416667, 0.0833334, 502.1450692400854, 0.3, 0.4, 0.11, 0.37)
)

$startRow = 12
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Reuse the same style applied to the other rows' column A cells
    # (bold font + border + centered alignment) instead of creating a new style.
    $ws.Range("A2").Copy($ws.Cells.Item($r, 1))

    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
